# realised that the peaks are 270 and 180GHz. Updating graphs for new
# experimental findings. Append three new rows (28-30) of interferogram
# mock-data parameters to the "1D NEW" sheet's table (Table43), expanding
# the table/selection to match, and note the realistic-interferogram
# observation in a quote-prefixed comment cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1D NEW")

# Row 31: 1dmockanderrors28.csv
$ws.Range("B31").Value = "1dmockanderrors28.csv"
$ws.Range("C31").Value = 99
$ws.Range("D31").Value = 500
$ws.Range("E31").Value = 0.3
$ws.Range("F31").Value = 0.05
$ws.Range("G31").Value = 360
$ws.Range("H31").Value = 1
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 5
$ws.Range("L31").Value = 1

# Row 32: 1dmockanderrors29.csv (noiseless)
$ws.Range("B32").Value = "1dmockanderrors29.csv"
$ws.Range("C32").Value = 99
$ws.Range("D32").Value = 500
$ws.Range("E32").Value = 0.3
$ws.Range("F32").Value = 0.05
$ws.Range("G32").Value = 360
$ws.Range("H32").Value = 1
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 1
$ws.Range("M32").Value = "noiseless"

# Row 33: 1dmockanderrors30.csv (realistic interferogram - peaks at 0.182 & 0.273 THz)
$ws.Range("B33").Value = "1dmockanderrors30.csv"
$ws.Range("C33").Value = 99
$ws.Range("D33").Value = 500
$ws.Range("E33").Value = "[Two equal peaks at 0.182 and 0.273]"
$ws.Range("F33").Value = "[modulated by a gaussian with FWHM = 35mm]"
$ws.Range("G33").Value = 360
$ws.Range("H33").Value = 1
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 1
# Leading apostrophe -> text forced/quote-prefixed, matching the source cell style
$ws.Range("M33").Value = "'realistic interferogram, very simular to what was observed"

# Grow the table (and its AutoFilter) to cover the three new rows
$tbl = $ws.ListObjects.Item("Table43")
$tbl.Resize($ws.Range("B3:M33"))

# Move the active selection down past the new data, like a user would after typing
$ws.Range("M34").Select()
